$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "compounds" sheet ---
$compounds = $wb.Worksheets.Item("compounds")
$compounds.Activate() | Out-Null
$compounds.Range("B1:G1").Select() | Out-Null

# --- Add the new "components" sheet right after "compounds" ---
$components = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $compounds)
$components.Name = "components"

# Header row
$components.Range("B1").Value = "as-is"
$components.Range("C1").Value = "almost as-is"
$components.Range("D1").Value = "implement"
$components.Range("E1").Value = "not supported"
$components.Range("F1").Value = "ionize"
$components.Range("G1").Value = "done"

# Row 2: $
$components.Range("A2").Value = "`$"
$components.Range("B2").Value = "X"
$components.Range("G2").Value = "X"

# Row 3: [
$components.Range("A3").Value = "["
$components.Range("C3").Value = "X"
$components.Range("G3").Value = "X"

# Row 4: [[
$components.Range("A4").Value = "[["
$components.Range("B4").Value = "X"
$components.Range("G4").Value = "X"

# Row 5: as.data.table
$components.Range("A5").Value = "as.data.table"
$components.Range("B5").Value = "X"
$components.Range("G5").Value = "X"

# Row 6: componentInfo
$components.Range("A6").Value = "componentInfo"
$components.Range("B6").Value = "X"
$components.Range("G6").Value = "X"

# Row 7: componentTable
$components.Range("A7").Value = "componentTable"
$components.Range("B7").Value = "X"
$components.Range("G7").Value = "X"

# Row 8: consensus
$components.Range("A8").Value = "consensus"
$components.Range("E8").Value = "X"

# Row 9: filter
$components.Range("A9").Value = "filter"
$components.Range("C9").Value = "X"
$components.Range("G9").Value = "X"

# Row 10: findFGroup
$components.Range("A10").Value = "findFGroup"
$components.Range("B10").Value = "X"
$components.Range("G10").Value = "X"

# Row 11: groupNames
$components.Range("A11").Value = "groupNames"
$components.Range("B11").Value = "X"
$components.Range("G11").Value = "X"

# Row 12: initialize
$components.Range("A12").Value = "initialize"
$components.Range("C12").Value = "X"
$components.Range("G12").Value = "X"

# Row 13: length
$components.Range("A13").Value = "length"
$components.Range("B13").Value = "X"
$components.Range("G13").Value = "X"

# Row 14: names
$components.Range("A14").Value = "names"
$components.Range("B14").Value = "X"
$components.Range("G14").Value = "X"

# Row 15: plotEIC
$components.Range("A15").Value = "plotEIC"
$components.Range("B15").Value = "X"
$components.Range("D15").Value = "X"
$components.Range("G15").Value = "X"
$components.Range("H15").Value = "Seems enough, assuming we're not planning to merge components"

# Row 16: plotEICHash
$components.Range("A16").Value = "plotEICHash"
$components.Range("B16").Value = "X"
$components.Range("G16").Value = "X"

# Row 17: plotSpec
$components.Range("A17").Value = "plotSpec"
$components.Range("B17").Value = "X"
$components.Range("D17").Value = "X"
$components.Range("G17").Value = "X"
$components.Range("H17").Value = "Seems enough, assuming we're not planning to merge components"

# Row 18: plotSpecHash
$components.Range("A18").Value = "plotSpecHash"
$components.Range("B18").Value = "X"
$components.Range("G18").Value = "X"

# Row 19: show
$components.Range("A19").Value = "show"
$components.Range("C19").Value = "X"
$components.Range("G19").Value = "X"

# Column A width (approximate best-fit)
$components.Columns("A:A").ColumnWidth = 16

# Selection / active cell on the new sheet
$components.Range("H17").Select() | Out-Null
